$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '63.677.24'
Set-TextCell $ws.Range('E2') '  +2.74%  '
Set-TextCell $ws.Range('D3') '2.490.45'
Set-TextCell $ws.Range('E3') '  +2.98%  '
Set-TextCell $ws.Range('E4') '  +0.28%  '
Set-TextCell $ws.Range('D5') '575.29'
Set-TextCell $ws.Range('E5') '  +1.83%  '
Set-TextCell $ws.Range('D6') '149.61'
Set-TextCell $ws.Range('E6') '  +4.50%  '
Set-TextCell $ws.Range('D7') '0.999'
Set-TextCell $ws.Range('E7') '  -0.16%  '
Set-TextCell $ws.Range('D8') '0.541'
Set-TextCell $ws.Range('E8') '  +1.80%  '
Set-TextCell $ws.Range('E9') '  +4.47%  '
Set-TextCell $ws.Range('E10') '  +0.34%  '
Set-TextCell $ws.Range('B11') 'Toncoin'
Set-TextCell $ws.Range('C11') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws.Range('D11') '5.38'
Set-TextCell $ws.Range('E11') '  +3.09%  '
Set-TextCell $ws.Range('B12') 'Cardano'
Set-TextCell $ws.Range('C12') 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell $ws.Range('D12') '0.365'
Set-TextCell $ws.Range('E12') '  +4.28%  '
Set-TextCell $ws.Range('D13') '27.30'
Set-TextCell $ws.Range('E13') '  +5.98%  '
Set-TextCell $ws.Range('D14') '0.0000186'
Set-TextCell $ws.Range('E14') '  +6.80%  '
Set-TextCell $ws.Range('D15') '2.930.86'
Set-TextCell $ws.Range('E15') '  +2.69%  '
Set-TextCell $ws.Range('D16') '63.672.45'
Set-TextCell $ws.Range('E16') '  +2.77%  '
Set-TextCell $ws.Range('D17') '2.500.41'
Set-TextCell $ws.Range('E17') '  +3.23%  '
Set-TextCell $ws.Range('E18') '  +2.91%  '
Set-TextCell $ws.Range('D19') '7.33'
Set-TextCell $ws.Range('E19') '  +6.95%  '
Set-TextCell $ws.Range('B20') 'Polkadot'
Set-TextCell $ws.Range('C20') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws.Range('D20') '4.26'
Set-TextCell $ws.Range('E20') '  +2.85%  '
Set-TextCell $ws.Range('B21') 'BitcoinCash'
Set-TextCell $ws.Range('C21') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws.Range('D21') '329.44'
Set-TextCell $ws.Range('E21') '  +2.04%  '
Set-TextCell $ws.Range('D22') '1.01'
Set-TextCell $ws.Range('E22') '  +0.52%  '
Set-TextCell $ws.Range('D23') '1.90'
Set-TextCell $ws.Range('E23') '  +8.60%  '
Set-TextCell $ws.Range('D24') '67.58'
Set-TextCell $ws.Range('E24') '  +1.53%  '
Set-TextCell $ws.Range('D25') '648.46'
Set-TextCell $ws.Range('E25') '  +15.15%  '
Set-TextCell $ws.Range('D26') '8.89'
Set-TextCell $ws.Range('E26') '  +1.03%  '
Set-TextCell $ws.Range('D27') '0.0000106'
Set-TextCell $ws.Range('E27') '  +12.18%  '
Set-TextCell $ws.Range('B28') 'WrappedeETH'
Set-TextCell $ws.Range('C28') 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextCell $ws.Range('D28') '2.601.16'
Set-TextCell $ws.Range('E28') '  +2.51%  '
Set-TextCell $ws.Range('B29') 'Fetch.AI'
Set-TextCell $ws.Range('C29') 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws.Range('D29') '1.53'
Set-TextCell $ws.Range('E29') '  +8.92%  '
Set-TextCell $ws.Range('D30') '8.61'
Set-TextCell $ws.Range('E30') '  +4.95%  '
Set-TextCell $ws.Range('D31') '0.997'
Set-TextCell $ws.Range('E31') '  -0.50%  '
Set-TextCell $ws.Range('D32') '0.145'
Set-TextCell $ws.Range('E32') '  -1.63%  '
Set-TextCell $ws.Range('D33') '1.92'
Set-TextCell $ws.Range('E33') '  +2.69%  '
Set-TextCell $ws.Range('D34') '5.21'
Set-TextCell $ws.Range('E34') '  +9.13%  '
Set-TextCell $ws.Range('D35') '1.56'
Set-TextCell $ws.Range('E35') '  +3.40%  '
Set-TextCell $ws.Range('D36') '0.389'
Set-TextCell $ws.Range('E36') '  +2.30%  '
Set-TextCell $ws.Range('D37') '0.998'
Set-TextCell $ws.Range('E37') '  -0.16%  '
Set-TextCell $ws.Range('D38') '5.58'
Set-TextCell $ws.Range('E38') '  +2.57%  '
Set-TextCell $ws.Range('D39') '19.02'
Set-TextCell $ws.Range('E39') '  +2.63%  '
Set-TextCell $ws.Range('E40') '  +2.23%  '
Set-TextCell $ws.Range('D41') '148.07'
Set-TextCell $ws.Range('E41') '  -4.65%  '
Set-TextCell $ws.Range('D42') '2.66'
Set-TextCell $ws.Range('E42') '  +16.68%  '
Set-TextCell $ws.Range('E43') '  +0.22%  '
Set-TextCell $ws.Range('D44') '153.32'
Set-TextCell $ws.Range('E44') '  +3.62%  '
Set-TextCell $ws.Range('E45') '  +4.49%  '
Set-TextCell $ws.Range('D46') '21.29'
Set-TextCell $ws.Range('E46') '  +7.17%  '
Set-TextCell $ws.Range('D47') '0.0551'
Set-TextCell $ws.Range('E47') '  +4.30%  '
Set-TextCell $ws.Range('E48') '  +3.69%  '
Set-TextCell $ws.Range('E49') '  +5.47%  '
Set-TextCell $ws.Range('D50') '0.0931'
Set-TextCell $ws.Range('E50') '  +1.05%  '
Set-TextCell $ws.Range('B51') 'ONDO'
Set-TextCell $ws.Range('C51') 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell $ws.Range('D51') '0.759'
Set-TextCell $ws.Range('E51') '  +7.05%  '
